$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking prices that must remain stored as TEXT
# (matching the source data which uses inline strings, not numbers).
# Force text format before writing, then restore the default "Normal" style
# so no stray number-format/style is left behind on the cell.
$dCells = @("D2", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D40", "D41", "D42", "D43", "D44", "D45")
foreach ($cell in $dCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = "243.87"
$ws.Range("D4").Value = "5.414"
$ws.Range("D5").Value = "0.05981"
$ws.Range("D6").Value = "3.463"
$ws.Range("D7").Value = "6.523"
$ws.Range("D8").Value = "0.8132"
$ws.Range("D9").Value = "0.9143"
$ws.Range("D10").Value = "0.1408"
$ws.Range("D11").Value = "0.07415"
$ws.Range("D12").Value = "0.03237"
$ws.Range("D13").Value = "0.03088"
$ws.Range("D14").Value = "0.09346"
$ws.Range("D15").Value = "3.854"
$ws.Range("D16").Value = "0.001572"
$ws.Range("D17").Value = "0.04682"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "0.0005939"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").Value = "0.006062"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").Value = "0.005010"
$ws.Range("E20").Value = "19HotbitTokenHTBBestin24h"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").Value = "0.0009859"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").Value = "0.00007801"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "3.613"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").Value = "2.130"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("D25").Value = "0.3204"
$ws.Range("D26").Value = "0.1302"
$ws.Range("D40").Value = "0.03937"
$ws.Range("D41").Value = "0.006221"
$ws.Range("D42").Value = "0.1076"
$ws.Range("D43").Value = "0.002621"
$ws.Range("D44").Value = "0.006423"
$ws.Range("D45").Value = "0.00005223"

foreach ($cell in $dCells) {
    $ws.Range($cell).Style = "Normal"
}
